# edit.ps1 - Adds a new paragraph describing the first database query,
# along with a footnote about Copilot usage, right after the paragraph
# ending in "...fifth century." and before the trailing empty paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document's body text
# (the one ending in "...fifth century. "), which is the second paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*fifth century.*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not locate the paragraph ending in 'fifth century.'"
}

# Collapse to the end of that paragraph and insert a brand-new paragraph
# after it; Word copies the paragraph formatting (spacing/indent/fonts)
# from the source paragraph automatically.
$rng = $target.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
# $rng now sits collapsed at the very start of the freshly created paragraph.

# --- Insert the main run of new text (first three sentences/segments) ---
$rng.InsertAfter("My first query is admittedly not as exciting as I think querying this database could be, but the database itself was such an undertaking to create that I chose to keep it simple for my initial exploration to best feel out the structure within the timeline of this assignment. I am admittedly not very comfortable with SQL, especially in the Python environment, and am still wary of breaking the database structure. The intent of the hagiographical database is to demonstrate the archetype of religious exchange, but I noticed an interesting trend that I thought could be effectively queried using the Saints information data. Almost every saint was born somewhere in Ireland, but a regional pattern was beginning to emerge")
$rng.Collapse(0)

$rng.InsertAfter(". Leinster and Munster, two provinces in Ireland, appeared as frequent birthplaces for saints. I wanted to pull them out of the table and analyze the percentage of saints in the database with birthplaces listed specifically as Leinster or Munster. ")
$rng.Collapse(0)

$rng.InsertAfter("I wrote this as two queries, so that I could arrange them by sex.")
$rng.Collapse(0)

# Remember this position: the footnote reference mark must be inserted
# here, i.e. right before the trailing explanatory sentence.
$footnoteMark = $rng.Start

# --- Insert the remaining sentence that follows the footnote reference ---
$rng.InsertAfter(" Only one female saint was born in either place. 3 male saints were born in Leinster, and 2 were born in Munster. This provides an avenue for further research: Why are those locations significantly overrepresented by male saints, and what is happening there to produce so many saints overall?")
$rng.Collapse(0)

# --- Now add the footnote reference at the remembered position ---
$fnPoint = $d.Range($footnoteMark, $footnoteMark)
$fn = $d.Footnotes.Add($fnPoint)

# Make sure the footnote reference mark carries the same Times New Roman
# font as the rest of the paragraph (in addition to the FootnoteReference
# character style it already has).
$fn.Reference.Font.NameAscii = "Times New Roman"
$fn.Reference.Font.NameOther = "Times New Roman"
$fn.Reference.Font.NameBi = "Times New Roman"

# --- Fill in the footnote body text ---
$fn.Range.Text = " I did leave copilot on while writing my queries. It was showing suggestions, but before I accepted any, mostly just a reminder to print the query, I double checked in my notes. I have cited copilot in the works cited."
